$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates: bulletin number 4 -> 5, week dates shift by one week ---
$ws.Range("A8").Value = "Volume 33   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/26/2026  Through  2/1/2026"

# --- Cells that flip between a numeric value and the text placeholder "0" ---
# must copy a same-row cell that already carries the desired style + type,
# since a plain .Value assignment cannot change the stored type/style.
$ws.Range("D15").Copy($ws.Range("C15"))    # -> text "0" (was numeric 2)
$ws.Range("D27").Copy($ws.Range("C27"))    # -> text "0" (was numeric 2)
$ws.Range("F28").Copy($ws.Range("C28"))    # -> numeric (was text "0")
$ws.Range("C28").Value = 2

# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = 9.090909090909
$ws.Range("L16").Value = -45.454545454545
$ws.Range("M16").Value = -52
$ws.Range("N16").Value = -91.111111111111
# Row 17
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 18
$ws.Range("K17").Value = 11.111111111111
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 81.818181818181
$ws.Range("N17").Value = -25.925925925925
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -36.842105263157
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = -27.272727272727
$ws.Range("L18").Value = -15.789473684210
$ws.Range("M18").Value = -23.809523809523
$ws.Range("N18").Value = -87.692307692307
# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 52
$ws.Range("H19").Value = 40.540540540540
$ws.Range("I19").Value = 58
$ws.Range("J19").Value = 39
$ws.Range("K19").Value = 48.717948717948
$ws.Range("L19").Value = -20.547945205479
$ws.Range("M19").Value = 61.111111111111
$ws.Range("N19").Value = -20.547945205479
# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 233.333333333333
$ws.Range("I20").Value = 34
$ws.Range("J20").Value = 11
$ws.Range("K20").Value = 209.090909090909
$ws.Range("L20").Value = 240
$ws.Range("M20").Value = 88.888888888888
$ws.Range("N20").Value = -83
# Row 21
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 26.315789473684
$ws.Range("F21").Value = 127
$ws.Range("G21").Value = 92
$ws.Range("H21").Value = 38.043478260869
$ws.Range("I21").Value = 144
$ws.Range("J21").Value = 103
$ws.Range("K21").Value = 39.805825242718
$ws.Range("L21").Value = 2.857142857142
$ws.Range("M21").Value = 28.571428571428
$ws.Range("N21").Value = -74.603174603174
# Row 22
$ws.Range("C22").Value = 2
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 20
# Row 24
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -8.695652173913
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 129
$ws.Range("H24").Value = -7.751937984496
$ws.Range("I24").Value = 125
$ws.Range("J24").Value = 149
$ws.Range("K24").Value = -16.107382550335
$ws.Range("L24").Value = -30.167597765363
$ws.Range("M24").Value = 127.272727272727
# Row 25
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 76
$ws.Range("H25").Value = -15.789473684210
$ws.Range("I25").Value = 67
$ws.Range("J25").Value = 91
$ws.Range("K25").Value = -26.373626373626
$ws.Range("L25").Value = -38.532110091743
# Row 26
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = 21.875
$ws.Range("I26").Value = 45
$ws.Range("J26").Value = 38
$ws.Range("K26").Value = 18.421052631578
$ws.Range("L26").Value = 21.621621621621
$ws.Range("M26").Value = 7.142857142857
# Row 27
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 300
# Row 28
$ws.Range("I28").Value = 3
$ws.Range("K28").Value = 200
$ws.Range("L28").Value = -50
